# Auto-generated edit script applying the Adamantoise_Profits numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 43148760
$ws.Range("I86").Value = 12473.685
$ws.Range("K86").Value = 12473.685
$ws.Range("M86").Value = -11350.685

$ws.Range("H89").Value = 43148760
$ws.Range("I89").Value = 12473.685
$ws.Range("K89").Value = 62368.425
$ws.Range("M89").Value = -56752.425

$ws.Range("H116").Value = 100037496
$ws.Range("J116").Value = 5499
$ws.Range("L116").Value = 5499
$ws.Range("N116").Value = -12383

$ws.Range("H132").Value = 1470.8
$ws.Range("I132").Value = 1504.2858
$ws.Range("K132").Value = 4512.857400000001
$ws.Range("M132").Value = -1982.857400000001

$ws.Range("H138").Value = 2836.6396
$ws.Range("I138").Value = 1827.6207
$ws.Range("K138").Value = 5482.8621
$ws.Range("M138").Value = -342.8621000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1906.44
$ws.Range("I2").Value = 1934.4783
$ws.Range("J2").Value = 1584
$ws.Range("K2").Value = 1934.4783
$ws.Range("L2").Value = 1584
$ws.Range("M2").Value = -1821.4783
$ws.Range("N2").Value = -1810

$ws.Range("H24").Value = 100355
$ws.Range("J24").Value = 100355
$ws.Range("L24").Value = 100355
$ws.Range("N24").Value = -101103

$ws.Range("H28").Value = 16333.333
$ws.Range("I28").Value = 16333.333
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 16333.333
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -16141.333
$ws.Range("N28").ClearContents()

$ws.Range("H61").Value = 3480.2727
$ws.Range("I61").Value = 3285.2222
$ws.Range("K61").Value = 3285.2222
$ws.Range("M61").Value = -3073.2222

$ws.Range("H92").Value = 67499.5
$ws.Range("J92").Value = 67499.5
$ws.Range("L92").Value = 67499.5
$ws.Range("N92").Value = -72491.5

$ws.Range("H93").Value = 100716.5
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 100716.5
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 100716.5
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -105708.5

$ws.Range("H94").Value = 79973.5
$ws.Range("J94").Value = 79973.5
$ws.Range("L94").Value = 79973.5
$ws.Range("N94").Value = -81775.5

$ws.Range("H96").Value = 80550.664
$ws.Range("J96").Value = 80550.664
$ws.Range("L96").Value = 80550.664
$ws.Range("N96").Value = -86042.664

$ws.Range("H97").Value = 953.44446
$ws.Range("J97").Value = 1380
$ws.Range("L97").Value = 1380
$ws.Range("N97").Value = -2372

$ws.Range("H98").Value = 70355
$ws.Range("J98").Value = 70355
$ws.Range("L98").Value = 70355
$ws.Range("N98").Value = -76345

$ws.Range("H99").Value = 16333.333
$ws.Range("I99").Value = 16333.333
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 16333.333
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -13338.333
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 100355
$ws.Range("J100").Value = 100355
$ws.Range("L100").Value = 100355
$ws.Range("N100").Value = -102519

$ws.Range("H101").Value = 53744.5
$ws.Range("J101").Value = 53744.5
$ws.Range("L101").Value = 53744.5
$ws.Range("N101").Value = -60234.5

$ws.Range("H102").Value = 2925.625
$ws.Range("I102").Value = 1181
$ws.Range("K102").Value = 1181
$ws.Range("M102").Value = 441

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H104").Value = 100612.5
$ws.Range("J104").Value = 100612.5
$ws.Range("L104").Value = 100612.5
$ws.Range("N104").Value = -107600.5

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H116").Value = 1906.44
$ws.Range("I116").Value = 1934.4783
$ws.Range("J116").Value = 1584
$ws.Range("K116").Value = 1934.4783
$ws.Range("L116").Value = 1584
$ws.Range("M116").Value = 359.5217
$ws.Range("N116").Value = -6172

$ws.Range("H136").Value = 3480.2727
$ws.Range("I136").Value = 3285.2222
$ws.Range("K136").Value = 9855.6666
$ws.Range("M136").Value = -7305.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1906.44
$ws.Range("I3").Value = 1934.4783
$ws.Range("J3").Value = 1584
$ws.Range("K3").Value = 1934.4783
$ws.Range("L3").Value = 1584
$ws.Range("M3").Value = -1820.4783
$ws.Range("N3").Value = -1812

$ws.Range("H134").Value = 4468971
$ws.Range("I134").Value = 6496620.5
$ws.Range("J134").Value = 8142.4
$ws.Range("K134").Value = 19489861.5
$ws.Range("L134").Value = 24427.2
$ws.Range("M134").Value = -19487326.5
$ws.Range("N134").Value = -29497.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 677.7143
$ws.Range("I22").Value = 657.3333
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 657.3333
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -307.3333
$ws.Range("N22").Value = -1500

$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H31").Value = 6840.7417
$ws.Range("I31").Value = 1557.7333
$ws.Range("K31").Value = 1557.7333
$ws.Range("M31").Value = -1262.7333

$ws.Range("H34").Value = 6840.7417
$ws.Range("I34").Value = 1557.7333
$ws.Range("K34").Value = 1557.7333
$ws.Range("M34").Value = -1355.7333

$ws.Range("H58").Value = 2543.228
$ws.Range("I58").Value = 2303.4695
$ws.Range("K58").Value = 2303.4695
$ws.Range("M58").Value = -2100.4695

$ws.Range("H132").Value = 3013.162
$ws.Range("I132").Value = 2674.8125
$ws.Range("K132").Value = 8024.4375
$ws.Range("M132").Value = -5494.4375

$ws.Range("H134").Value = 1986.5366
$ws.Range("I134").Value = 1701.6857
$ws.Range("K134").Value = 5105.0571
$ws.Range("M134").Value = -2570.0571

$ws.Range("H136").Value = 2543.228
$ws.Range("I136").Value = 2303.4695
$ws.Range("K136").Value = 6910.4085
$ws.Range("M136").Value = -4360.4085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1674.8
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H117").Value = 4859.364
$ws.Range("J117").Value = 4896.6
$ws.Range("L117").Value = 14689.8
$ws.Range("N117").Value = -21573.8

$ws.Range("H121").Value = 5039705.5
$ws.Range("I121").Value = 840.3333
$ws.Range("K121").Value = 2520.9999
$ws.Range("M121").Value = -1210.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7635.5454
$ws.Range("I132").Value = 7998.857
$ws.Range("K132").Value = 23996.571
$ws.Range("M132").Value = -21466.571

$ws.Range("H135").Value = 50000
$ws.Range("I135").Value = 50000
$ws.Range("K135").Value = 50000
$ws.Range("M135").Value = -44930

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2199.111
$ws.Range("I82").Value = 1974
$ws.Range("K82").Value = 1974
$ws.Range("M82").Value = -1613

$ws.Range("H85").Value = 2199.111
$ws.Range("I85").Value = 1974
$ws.Range("K85").Value = 1974
$ws.Range("M85").Value = -726

$ws.Range("H136").Value = 3210.889
$ws.Range("J136").Value = 4099.8
$ws.Range("L136").Value = 12299.4
$ws.Range("N136").Value = -17399.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 21349.25
$ws.Range("J56").Value = 21349.25
$ws.Range("L56").Value = 21349.25
$ws.Range("N56").Value = -22777.25

$ws.Range("H94").Value = 78500
$ws.Range("J94").Value = 78500
$ws.Range("L94").Value = 78500
$ws.Range("N94").Value = -80302

$ws.Range("H113").Value = 1002.5
$ws.Range("I113").Value = 1002
$ws.Range("K113").Value = 3006
$ws.Range("M113").Value = -836

$ws.Range("H136").Value = 2889.6365
$ws.Range("I136").Value = 2332.5386
$ws.Range("J136").Value = 3694.3333
$ws.Range("K136").Value = 6997.6158
$ws.Range("L136").Value = 11082.9999
$ws.Range("M136").Value = -4447.6158
$ws.Range("N136").Value = -16182.9999

